$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# 1) "2016-08-18 18:15:43" -> "2016-08-18 18:16:44"
#    Overview!G2, Overview!G3 ; de-de!H2, de-de!H3
$wsOverview.Range("G2").Value = "2016-08-18 18:16:44"
$wsOverview.Range("G3").Value = "2016-08-18 18:16:44"
$wsDeDe.Range("H2").Value = "2016-08-18 18:16:44"
$wsDeDe.Range("H3").Value = "2016-08-18 18:16:44"

# 2) "ht" -> "mt" (Priority column), shared across zh-cn and de-de rows 2 & 3
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# 3) "2016-08-18 18:15:38" -> "2016-08-18 18:16:39" (zh-cn Correspond Handoff Datetime rows 2 & 3)
$wsZhCn.Range("H2").Value = "2016-08-18 18:16:39"
$wsZhCn.Range("H3").Value = "2016-08-18 18:16:39"

# 4) "2016-08-18 18:16:08" -> "2016-08-18 18:16:56" (zh-cn Correspond Handback DateTime rows 2 & 3)
$wsZhCn.Range("K2").Value = "2016-08-18 18:16:56"
$wsZhCn.Range("K3").Value = "2016-08-18 18:16:56"

# 5) "2016-08-18 18:16:17" -> "2016-08-18 18:17:10" (de-de Correspond Handback DateTime rows 2 & 3)
$wsDeDe.Range("K2").Value = "2016-08-18 18:17:10"
$wsDeDe.Range("K3").Value = "2016-08-18 18:17:10"

$wb.Save()
